# Fill in the "결과" (result) row for the 964회차 block and the full
# "965회차" second-line block, then move the active selection to K14,
# matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 964회차 block - "결과" row (row 14), columns B..H
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = 38
$ws.Range("G14").Value = 40
$ws.Range("H14").Value = 9

# 965회차 block - second "라인" group (rows 9-13), columns K..P
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 8
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 30
$ws.Range("O9").Value = 38
$ws.Range("P9").Value = 39

$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 37
$ws.Range("O10").Value = 35
$ws.Range("P10").Value = 44

$ws.Range("K11").Value = 12
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 11
$ws.Range("N11").Value = 37
$ws.Range("O11").Value = 38
$ws.Range("P11").Value = 41

$ws.Range("K12").Value = 12
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 8
$ws.Range("N12").Value = 36
$ws.Range("O12").Value = 38
$ws.Range("P12").Value = 43

$ws.Range("K13").Value = 6
$ws.Range("L13").Value = 9
$ws.Range("M13").Value = 21
$ws.Range("N13").Value = 25
$ws.Range("O13").Value = 26
$ws.Range("P13").Value = 43

# Move the selection like the author left it
$ws.Range("K14").Select()
